$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 483, shifting the existing rows (483:557) down to (484:558).
$ws.Rows(483).Insert()

# Populate the newly inserted row 483 with the new price entry.
$ws.Cells.Item(483, 1).Value = 10
$ws.Cells.Item(483, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(483, 3).Value = "La Araucanía"
$ws.Cells.Item(483, 4).Value = 44984
$ws.Cells.Item(483, 5).Value = 9
$ws.Cells.Item(483, 6).Value = 100112008
$ws.Cells.Item(483, 7).Value = "Coliflor"
$ws.Cells.Item(483, 8).Value = "Sin especificar"
$ws.Cells.Item(483, 9).Value = "Primera"
$ws.Cells.Item(483, 10).Value = 1000
$ws.Cells.Item(483, 11).Value = 1500
$ws.Cells.Item(483, 12).Value = 1500
$ws.Cells.Item(483, 13).Value = 1500
$ws.Cells.Item(483, 14).Value = "`$/unidad"
$ws.Cells.Item(483, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(483, 16).Value = 1500
$ws.Cells.Item(483, 17).Value = 1
$ws.Cells.Item(483, 18).Value = "Hortaliza"
